# "add orders list on user home page" - append a new hashed value in column C
# next to the last email/password row, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "6f42205f20666a04d1f2ee777a799c383d4af4be593eb607b0f1cda06bc73c50"

[void]$ws.Range("B8").Select()
